# Modifies assert for main nav: the MainNavigationBarTest block on the
# "Test Data" sheet now inserts a "What" category (between CategoryName and
# Necklaces) and renames "Guest Bartender" -> "Guest", pushing "Minibar"
# down into a new row. The MainNavSubCategoryTest block below is unchanged
# in content but shifts down by one row to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Data")

# Insert a new row above the old "MainNavSubCategoryTest" header row (56)
# so everything from there on shifts down by one.
$ws.Rows.Item(56).Insert()

# --- MainNavigationBarTest block (rows 46-55) ---
# Row 46 (header row "CategoryName") stays the same.
# Write "Guest" before "What" so the two brand-new shared-string entries
# land in the same order as the authored workbook (Guest=44, What=45).
$ws.Range("B54").Value = "Guest"
$ws.Range("B47").Value = "What"
$ws.Range("B48").Value = "Necklaces"
$ws.Range("B49").Value = "Bracelets"
$ws.Range("B50").Value = "Earrings"

# Row 51 used to be "Featured Shops" (plain style) and now becomes "Rings"
# using the same styling as rows 47-50 (text format + explicit black font)
# plus a (blank) C cell to match.
$ws.Range("B51").Value = "Rings"
$ws.Range("B51").Font.Color = 0
$ws.Range("C51").Value = $null
$ws.Range("C51").Font.Color = 0

$ws.Range("B52").Value = "Featured Shops"
$ws.Range("B53").Value = "Personalized"

# New row 55 holds what used to trail row 54 ("Minibar").
$ws.Range("A55").Value = "Y"
$ws.Range("A55").Font.Color = 0
$ws.Range("B55").Value = "Minibar"

# --- MainNavSubCategoryTest block (rows 57-67), shifted down by one ---
# Rows 57-66 already hold the former 56-65 content after the row insert
# above. Row 63 gains the same "blank C cell" treatment as row 51 did.
$ws.Range("B63").Font.Color = 0
$ws.Range("C63").Value = $null
$ws.Range("C63").Font.Color = 0

# Row 67 is the new last data row, mirroring what used to be row 66
# ("Minibar").
$ws.Range("A67").Value = "Y"
$ws.Range("A67").Font.Color = 0
$ws.Range("B67").Value = "Minibar"

# Fix the selection to match the authored state: active cell B47, single
# cell selected (rather than the whole row).
$ws.Range("B47").Select()
